$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final (post-edit) player table, in on-sheet row order (A2:C17).
# The underlying change simply re-orders the existing player rows;
# every name/position/team combination that existed before still
# exists after the edit, just shuffled to new rows.
$data = @(
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Jalen Johnson", "PF", "Atlanta Hawks"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Mark Williams", "C", "Charlotte Hornets"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

$wb.Save()
